$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new survey submission ("南京叶子科技有限公司（南京OPPO）") was received and
# inserted as the new first data row. All existing data rows (previously
# rows 2-13) shift down by one (to rows 3-14); the running index in column A
# stays a simple 0-based positional counter.
# ---------------------------------------------------------------------------

# Insert a new row before row 2, shifting existing rows 2-13 down to 3-14.
$ws.Rows.Item(2).Insert()

# The inserted row picks up stray formatting from the insert operation;
# clear it, then reapply the original "index column" style
# (bold/bordered/centered) to A2 only -- matching A3:A14, which keep that
# same style while columns B:Q carry no explicit cell style.
$ws.Range("A2:Q2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 2 content (only company name and update time were actually filled
# in by the submitter; the rest keeps the template's placeholder text).
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "南京叶子科技有限公司（南京OPPO）"
$ws.Range("C2").Value = "xx区"
$ws.Range("D2").Value = "xxx事业部"
$ws.Range("E2").Value = "Java"
$ws.Range("F2").Value = "9:00-18:30"
$ws.Range("G2").Value = "1.5h"
$ws.Range("H2").Value = "135 加班，24 正常；大小周等等"
$ws.Range("I2").Value = "基数 xxxx，比例 xx%"
$ws.Range("J2").Value = "13薪还是根据公司业绩提供，是否折扣，折扣比例。"
$ws.Range("K2").Value = "是否打折，比如 xx%。"
$ws.Range("L2").Value = "工位大小，环境，是否提供设备，设备型号种类。"
$ws.Range("M2").Value = "是否有入职就有，是否有前置条件才有。"
$ws.Range("N2").Value = "是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。"
# O2 and P2 stay blank (left untouched after the row insert/clear above, so
# they remain present-but-empty cells, matching the template row's blanks).
$ws.Range("Q2").Value = "2022-01-25 03:43:17"

# Renumber the index column (A) for all data rows so it stays a simple
# positional counter (0, 1, 2, ...) regardless of content shifting.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Two cells hold numeric-looking text ("8%" and "995") that must stay text,
# not be auto-converted to a percentage/number. Force a text format before
# assigning, then clear the format again so no stray style index remains.
$ws.Cells.Item(6, 9).NumberFormat = "@"
$ws.Cells.Item(6, 9).Value = "8%"
$ws.Cells.Item(6, 9).ClearFormats()

$ws.Cells.Item(10, 8).NumberFormat = "@"
$ws.Cells.Item(10, 8).Value = "995"
$ws.Cells.Item(10, 8).ClearFormats()
